# Update "想去人数" (interest count) figures in column F for several rows
# across the 展览 (Exhibitions), 演出 (Shows) and 全部类型 (All types) sheets,
# matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1883
$ws1.Range("F6").Value  = 853
$ws1.Range("F12").Value = 16
$ws1.Range("F14").Value = 140
$ws1.Range("F16").Value = 4434
$ws1.Range("F20").Value = 430
$ws1.Range("F23").Value = 1110
$ws1.Range("F24").Value = 2061
$ws1.Range("F26").Value = 51
$ws1.Range("F27").Value = 30
$ws1.Range("F29").Value = 2128

# --- Sheet: 演出 (Shows) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 34

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1883
$ws4.Range("F6").Value  = 853
$ws4.Range("F12").Value = 16
$ws4.Range("F14").Value = 140
$ws4.Range("F16").Value = 34
$ws4.Range("F17").Value = 4434
$ws4.Range("F21").Value = 430
$ws4.Range("F24").Value = 1110
$ws4.Range("F25").Value = 2061
$ws4.Range("F27").Value = 51
$ws4.Range("F28").Value = 30
$ws4.Range("F30").Value = 2128
